# DataTable.xlsx update: refresh "Job adverts by occupation" latest period,
# and update the sheet's active selection as left by the author on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Job adverts by occupation" (row 13) - Latest period (release date) column C
# changes from "Nov 2024 (07/02/25)" to "Feb 2025 (01/04/25)"
$ws.Range("C13").Value = "Feb 2025 (01/04/25)"

# Reflect the scroll position / selection state saved with the workbook
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("C12").Select()
